$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 41.428665
$ws.Range("H2").Value = 124.285995
$ws.Range("I2").Value = 0.06969137269740189
$ws.Range("J2").Value = 0.06969137269740189
$ws.Range("M2").Value = 0.3331066666666667
$ws.Range("N2").Value = 0.99932
$ws.Range("O2").Value = 0.002125805913843485
$ws.Range("P2").Value = 0.002125805913843485
$ws.Range("Q2").Value = 13.8001645026
$ws.Range("R2").Value = 124.2014805234
$ws.Range("S2").Value = 0.0001481503322240073
$ws.Range("T2").Value = 0.0001481503322240073
# Row 3
$ws.Range("G3").Value = 41.428665
$ws.Range("H3").Value = 124.285995
$ws.Range("I3").Value = 0.06969137269740189
$ws.Range("J3").Value = 0.06969137269740189
$ws.Range("O3").Value = 0.0008775937418887864
$ws.Range("P3").Value = 0.0008775937418887864
$ws.Range("Q3").Value = 5.697104296139999
$ws.Range("R3").Value = 51.27393866526
$ws.Range("S3").Value = 0.00006116071254287892
$ws.Range("T3").Value = 0.00006116071254287892
# Row 4
$ws.Range("G4").Value = 41.428665
$ws.Range("H4").Value = 124.285995
$ws.Range("I4").Value = 0.06969137269740189
$ws.Range("J4").Value = 0.06969137269740189
$ws.Range("M4").Value = 91.40156066666667
$ws.Range("N4").Value = 274.204682
$ws.Range("O4").Value = 0.5833025803538128
$ws.Range("P4").Value = 0.5833025803538128
$ws.Range("Q4").Value = 3786.64463733651
$ws.Range("R4").Value = 34079.80173602859
$ws.Range("S4").Value = 0.04065115752279378
$ws.Range("T4").Value = 0.04065115752279378
# Row 5
$ws.Range("G5").Value = 41.428665
$ws.Range("H5").Value = 124.285995
$ws.Range("I5").Value = 0.06969137269740189
$ws.Range("J5").Value = 0.06969137269740189
$ws.Range("M5").Value = 0.5759770000000001
$ws.Range("N5").Value = 1.727931
$ws.Range("O5").Value = 0.00367574544541637
$ws.Range("P5").Value = 0.00367574544541637
$ws.Range("Q5").Value = 23.861958180705
$ws.Range("R5").Value = 214.757623626345
$ws.Range("S5").Value = 0.0002561677457772898
$ws.Range("T5").Value = 0.0002561677457772897
# Row 6
$ws.Range("G6").Value = 41.428665
$ws.Range("H6").Value = 124.285995
$ws.Range("I6").Value = 0.06969137269740189
$ws.Range("J6").Value = 0.06969137269740189
$ws.Range("M6").Value = 64.24849033333334
$ws.Range("N6").Value = 192.745471
$ws.Range("O6").Value = 0.4100182745450386
$ws.Range("P6").Value = 0.4100182745450385
$ws.Range("Q6").Value = 2661.729182775405
$ws.Range("R6").Value = 23955.56264497864
$ws.Range("S6").Value = 0.02857473638406393
$ws.Range("T6").Value = 0.02857473638406393
# Row 7
$ws.Range("H7").Value = 510.696747
$ws.Range("I7").Value = 0.2863649869040173
$ws.Range("J7").Value = 0.2863649869040173
$ws.Range("M7").Value = 0.3331066666666667
$ws.Range("N7").Value = 0.99932
$ws.Range("O7").Value = 0.002125805913843485
$ws.Range("P7").Value = 0.002125805913843485
$ws.Range("Q7").Value = 56.70549702356
$ws.Range("R7").Value = 510.3494732120399
$ws.Range("S7").Value = 0.0006087563826782721
$ws.Range("T7").Value = 0.0006087563826782722
# Row 8
$ws.Range("H8").Value = 510.696747
$ws.Range("I8").Value = 0.2863649869040173
$ws.Range("J8").Value = 0.2863649869040173
$ws.Range("O8").Value = 0.0008775937418887864
$ws.Range("P8").Value = 0.0008775937418887864
$ws.Range("S8").Value = 0.0002513121204030298
$ws.Range("T8").Value = 0.0002513121204030299
# Row 9
$ws.Range("H9").Value = 510.696747
$ws.Range("I9").Value = 0.2863649869040173
$ws.Range("J9").Value = 0.2863649869040173
$ws.Range("M9").Value = 91.40156066666667
$ws.Range("N9").Value = 274.204682
$ws.Range("O9").Value = 0.5833025803538128
$ws.Range("P9").Value = 0.5833025803538128
$ws.Range("Q9").Value = 15559.49323439661
$ws.Range("R9").Value = 140035.4391095694
$ws.Range("S9").Value = 0.1670374357840991
$ws.Range("T9").Value = 0.1670374357840991
# Row 10
$ws.Range("H10").Value = 510.696747
$ws.Range("I10").Value = 0.2863649869040173
$ws.Range("J10").Value = 0.2863649869040173
$ws.Range("M10").Value = 0.5759770000000001
$ws.Range("N10").Value = 1.727931
$ws.Range("O10").Value = 0.00367574544541637
$ws.Range("P10").Value = 0.00367574544541637
$ws.Range("Q10").Value = 98.04986008227301
$ws.Range("R10").Value = 882.448740740457
$ws.Range("S10").Value = 0.00105260479633916
$ws.Range("T10").Value = 0.00105260479633916
# Row 11
$ws.Range("H11").Value = 510.696747
$ws.Range("I11").Value = 0.2863649869040173
$ws.Range("J11").Value = 0.2863649869040173
$ws.Range("M11").Value = 64.24849033333334
$ws.Range("N11").Value = 192.745471
$ws.Range("O11").Value = 0.4100182745450386
$ws.Range("P11").Value = 0.4100182745450385
$ws.Range("Q11").Value = 10937.16500429809
$ws.Range("R11").Value = 98434.48503868283
$ws.Range("S11").Value = 0.1174148778204977
$ws.Range("T11").Value = 0.1174148778204978
# Row 12
$ws.Range("G12").Value = 244.5761666666666
$ws.Range("H12").Value = 733.7284999999999
$ws.Range("I12").Value = 0.4114264551867299
$ws.Range("J12").Value = 0.41142645518673
$ws.Range("M12").Value = 0.3331066666666667
$ws.Range("N12").Value = 0.99932
$ws.Range("O12").Value = 0.002125805913843485
$ws.Range("P12").Value = 0.002125805913843485
$ws.Range("Q12").Value = 81.46995162444443
$ws.Range("R12").Value = 733.2295646199999
$ws.Range("S12").Value = 0.000874612791547612
$ws.Range("T12").Value = 0.0008746127915476122
# Row 13
$ws.Range("G13").Value = 244.5761666666666
$ws.Range("H13").Value = 733.7284999999999
$ws.Range("I13").Value = 0.4114264551867299
$ws.Range("J13").Value = 0.41142645518673
$ws.Range("O13").Value = 0.0008775937418887864
$ws.Range("P13").Value = 0.0008775937418887864
$ws.Range("Q13").Value = 33.63313613533333
$ws.Range("R13").Value = 302.698225218
$ws.Range("S13").Value = 0.0003610652823193614
$ws.Range("T13").Value = 0.0003610652823193615
# Row 14
$ws.Range("G14").Value = 244.5761666666666
$ws.Range("H14").Value = 733.7284999999999
$ws.Range("I14").Value = 0.4114264551867299
$ws.Range("J14").Value = 0.41142645518673
$ws.Range("M14").Value = 91.40156066666667
$ws.Range("N14").Value = 274.204682
$ws.Range("O14").Value = 0.5833025803538128
$ws.Range("P14").Value = 0.5833025803538128
$ws.Range("Q14").Value = 22354.64333520411
$ws.Range("R14").Value = 201191.790016837
$ws.Range("S14").Value = 0.2399861129362419
$ws.Range("T14").Value = 0.2399861129362419
# Row 15
$ws.Range("G15").Value = 244.5761666666666
$ws.Range("H15").Value = 733.7284999999999
$ws.Range("I15").Value = 0.4114264551867299
$ws.Range("J15").Value = 0.41142645518673
$ws.Range("M15").Value = 0.5759770000000001
$ws.Range("N15").Value = 1.727931
$ws.Range("O15").Value = 0.00367574544541637
$ws.Range("P15").Value = 0.00367574544541637
$ws.Range("Q15").Value = 140.8702467481667
$ws.Range("R15").Value = 1267.8322207335
$ws.Range("S15").Value = 0.001512298918776425
$ws.Range("T15").Value = 0.001512298918776425
# Row 16
$ws.Range("G16").Value = 244.5761666666666
$ws.Range("H16").Value = 733.7284999999999
$ws.Range("I16").Value = 0.4114264551867299
$ws.Range("J16").Value = 0.41142645518673
$ws.Range("M16").Value = 64.24849033333334
$ws.Range("N16").Value = 192.745471
$ws.Range("O16").Value = 0.4100182745450386
$ws.Range("P16").Value = 0.4100182745450385
$ws.Range("Q16").Value = 15713.64947984705
$ws.Range("R16").Value = 141422.8453186235
$ws.Range("S16").Value = 0.1686923652578446
$ws.Range("T16").Value = 0.1686923652578446
# Row 17
$ws.Range("G17").Value = 24.173247
$ws.Range("H17").Value = 72.51974100000001
$ws.Range("I17").Value = 0.04066427836821081
$ws.Range("J17").Value = 0.04066427836821081
$ws.Range("M17").Value = 0.3331066666666667
$ws.Range("N17").Value = 0.99932
$ws.Range("O17").Value = 0.002125805913843485
$ws.Range("P17").Value = 0.002125805913843485
$ws.Range("Q17").Value = 8.052269730680001
$ws.Range("R17").Value = 72.47042757612002
$ws.Range("S17").Value = 0.00008644436343732025
$ws.Range("T17").Value = 0.00008644436343732025
# Row 18
$ws.Range("G18").Value = 24.173247
$ws.Range("H18").Value = 72.51974100000001
$ws.Range("I18").Value = 0.04066427836821081
$ws.Range("J18").Value = 0.04066427836821081
$ws.Range("O18").Value = 0.0008775937418887864
$ws.Range("P18").Value = 0.0008775937418887864
$ws.Range("Q18").Value = 3.324208234452001
$ws.Range("R18").Value = 29.91787411006801
$ws.Range("S18").Value = 0.00003568671621436536
$ws.Range("T18").Value = 0.00003568671621436536
# Row 19
$ws.Range("G19").Value = 24.173247
$ws.Range("H19").Value = 72.51974100000001
$ws.Range("I19").Value = 0.04066427836821081
$ws.Range("J19").Value = 0.04066427836821081
$ws.Range("M19").Value = 91.40156066666667
$ws.Range("N19").Value = 274.204682
$ws.Range("O19").Value = 0.5833025803538128
$ws.Range("P19").Value = 0.5833025803538128
$ws.Range("Q19").Value = 2209.472502180819
$ws.Range("R19").Value = 19885.25251962736
$ws.Range("S19").Value = 0.0237195785004031
$ws.Range("T19").Value = 0.0237195785004031
# Row 20
$ws.Range("G20").Value = 24.173247
$ws.Range("H20").Value = 72.51974100000001
$ws.Range("I20").Value = 0.04066427836821081
$ws.Range("J20").Value = 0.04066427836821081
$ws.Range("M20").Value = 0.5759770000000001
$ws.Range("N20").Value = 1.727931
$ws.Range("O20").Value = 0.00367574544541637
$ws.Range("P20").Value = 0.00367574544541637
$ws.Range("Q20").Value = 13.923234287319
$ws.Range("R20").Value = 125.309108585871
$ws.Range("S20").Value = 0.0001494715360030943
$ws.Range("T20").Value = 0.0001494715360030943
# Row 21
$ws.Range("G21").Value = 24.173247
$ws.Range("H21").Value = 72.51974100000001
$ws.Range("I21").Value = 0.04066427836821081
$ws.Range("J21").Value = 0.04066427836821081
$ws.Range("M21").Value = 64.24849033333334
$ws.Range("N21").Value = 192.745471
$ws.Range("O21").Value = 0.4100182745450386
$ws.Range("P21").Value = 0.4100182745450385
$ws.Range("Q21").Value = 1553.094626204779
$ws.Range("R21").Value = 13977.85163584301
$ws.Range("S21").Value = 0.01667309725215294
$ws.Range("T21").Value = 0.01667309725215293
# Row 22
$ws.Range("G22").Value = 114.0486906666667
$ws.Range("H22").Value = 342.146072
$ws.Range("I22").Value = 0.19185290684364
$ws.Range("J22").Value = 0.19185290684364
$ws.Range("M22").Value = 0.3331066666666667
$ws.Range("N22").Value = 0.99932
$ws.Range("O22").Value = 0.002125805913843485
$ws.Range("P22").Value = 0.002125805913843485
$ws.Range("Q22").Value = 37.99037918567111
$ws.Range("R22").Value = 341.91341267104
$ws.Range("S22").Value = 0.0004078420439562731
$ws.Range("T22").Value = 0.0004078420439562731
# Row 23
$ws.Range("G23").Value = 114.0486906666667
$ws.Range("H23").Value = 342.146072
$ws.Range("I23").Value = 0.19185290684364
$ws.Range("J23").Value = 0.19185290684364
$ws.Range("O23").Value = 0.0008775937418887864
$ws.Range("P23").Value = 0.0008775937418887864
$ws.Range("Q23").Value = 15.68351974571733
$ws.Range("R23").Value = 141.151677711456
$ws.Range("S23").Value = 0.0001683689104091508
$ws.Range("T23").Value = 0.0001683689104091508
# Row 24
$ws.Range("G24").Value = 114.0486906666667
$ws.Range("H24").Value = 342.146072
$ws.Range("I24").Value = 0.19185290684364
$ws.Range("J24").Value = 0.19185290684364
$ws.Range("M24").Value = 91.40156066666667
$ws.Range("N24").Value = 274.204682
$ws.Range("O24").Value = 0.5833025803538128
$ws.Range("P24").Value = 0.5833025803538128
$ws.Range("Q24").Value = 10424.22831892323
$ws.Range("R24").Value = 93818.0548703091
$ws.Range("S24").Value = 0.1119082956102749
$ws.Range("T24").Value = 0.1119082956102749
# Row 25
$ws.Range("G25").Value = 114.0486906666667
$ws.Range("H25").Value = 342.146072
$ws.Range("I25").Value = 0.19185290684364
$ws.Range("J25").Value = 0.19185290684364
$ws.Range("M25").Value = 0.5759770000000001
$ws.Range("N25").Value = 1.727931
$ws.Range("O25").Value = 0.00367574544541637
$ws.Range("P25").Value = 0.00367574544541637
$ws.Range("Q25").Value = 65.68942270411468
$ws.Range("R25").Value = 591.204804337032
$ws.Range("S25").Value = 0.0007052024485204009
$ws.Range("T25").Value = 0.0007052024485204008
# Row 26
$ws.Range("G26").Value = 114.0486906666667
$ws.Range("H26").Value = 342.146072
$ws.Range("I26").Value = 0.19185290684364
$ws.Range("J26").Value = 0.19185290684364
$ws.Range("M26").Value = 64.24849033333334
$ws.Range("N26").Value = 192.745471
$ws.Range("O26").Value = 0.4100182745450386
$ws.Range("P26").Value = 0.4100182745450385
$ws.Range("Q26").Value = 7327.456199826657
$ws.Range("R26").Value = 65947.10579843992
$ws.Range("S26").Value = 0.07866319783047929
$ws.Range("T26").Value = 0.07866319783047927
